$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TastingAnalysis")

# Fill in the previously-blank "April" row (row 28) with actuals.
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = 12
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 12
$ws.Range("I28").Value = 1
$ws.Range("K28").Value = 111111
$ws.Range("L28").Value = 123123
$ws.Range("M28").Value = 12312
$ws.Range("N28").Value = 123123
$ws.Range("O28").Formula = "=L28/H28"
$ws.Range("P28").Value = 11111

# Update the frozen-pane scroll position / active selection to match the
# author's last on-screen state.
$ws.Application.Goto($ws.Range("A17"), $false)
$ws.Range("D28").Select()
